$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 18:38"

# Country data table (country, casos totales, nuevos casos, casos activos, recuperados, casos criticos, muertes hoy, muertes)
$countryData = @(
    @("Estados Unidos", 8468223, 11570, 5513584, 2729188, 0, 229, 225451),
    @("India", 7644979, 50243, 6784742, 744358, 0, 643, 115879),
    @("Brasil", 5255277, 4150, 4681659, 419291, 0, 101, 154327),
    @("Rusia", 1431635, 16319, 1085608, 321392, 0, 269, 24635),
    @("España", 1015795, 0, 0, 0, 0, 0, 33992),
    @("Argentina", 1002662, 0, 803965, 171981, 0, 0, 26716),
    @("Colombia", 965883, 0, 867961, 68820, 0, 0, 29102),
    @("Francia", 910277, 0, 105935, 770719, 0, 0, 33623),
    @("Peru", 870876, 0, 784056, 53000, 0, 0, 33820),
    @("Mexico", 854926, 3699, 623494, 145094, 0, 171, 86338),
    @("Reino Unido", 762542, 21331, 0, 0, 0, 241, 43967),
    @("Sudafrica", 705254, 0, 635257, 51505, 0, 0, 18492),
    @("Iran", 539670, 5039, 434676, 73960, 0, 322, 31034),
    @("Chile", 494478, 1099, 466643, 14133, 0, 26, 13702),
    @("Irak", 434598, 3920, 366134, 58098, 0, 49, 10366),
    @("Italia", 434449, 10874, 255005, 142739, 0, 89, 36705),
    @("Banglades", 391586, 1380, 307141, 78746, 0, 18, 5699),
    @("Alemania", 376588, 2857, 298300, 68366, 0, 23, 9922),
    @("Indonesia", 368842, 3602, 293653, 62455, 0, 117, 12734),
    @("Filipinas", 360775, 1640, 310642, 43443, 0, 17, 6690),
    @("Turquia", 349519, 0, 305427, 34721, 0, 0, 9371),
    @("Arabia Saudita", 342968, 385, 329270, 8481, 0, 16, 5217),
    @("Pakistan", 324077, 618, 308020, 9384, 0, 14, 6673),
    @("Ucrania", 309107, 5469, 129533, 173788, 0, 113, 5786),
    @("Israel", 305633, 757, 279740, 23622, 0, 8, 2271),
    @("Paises Bajos", 244391, 8165, 0, 0, 0, 46, 6814),
    @("Belgica", 230480, 8227, 21214, 198823, 0, 30, 10443),
    @("Canada", 203046, 1609, 171239, 22015, 0, 14, 9792),
    @("Polonia", 192539, 9291, 95956, 92862, 0, 107, 3721),
    @("Chequia", 188773, 6811, 75645, 111528, 0, 87, 1600),
    @("Rumania", 186254, 3400, 134395, 45863, 0, 65, 5996),
    @("Marruecos", 175749, 0, 146421, 26352, 0, 0, 2976),
    @("Ecuador", 154115, 692, 134187, 7524, 0, 9, 12404),
    @("Bolivia", 139890, 119, 104957, 26431, 0, 21, 8502),
    @("Nepal", 139129, 3093, 96609, 41755, 0, 8, 765),
    @("Catar", 129671, 0, 126650, 2797, 0, 0, 224),
    @("Panama", 125181, 0, 101545, 21062, 0, 0, 2574),
    @("Republica Dominicana", 121973, 306, 99271, 20498, 0, 1, 2204),
    @("Kuwait", 117718, 886, 109198, 7806, 0, 4, 714),
    @("Emiratos Arabes Unidos", 117594, 1077, 110313, 6811, 0, 4, 470),
    @("Oman", 111033, 439, 96949, 12962, 0, 8, 1122),
    @("Kazajistan", 109623, 115, 105145, 2710, 0, 0, 1768),
    @("Suecia", 106380, 0, 0, 0, 0, 6, 5922),
    @("Egipto", 105547, 0, 98314, 1103, 0, 0, 6130),
    @("Portugal", 103736, 1876, 61898, 39625, 0, 15, 2213),
    @("Guatemala", 102219, 620, 91604, 7069, 0, 5, 3546),
    @("Costa Rica", 97075, 0, 59580, 36291, 0, 0, 1204),
    @("Japon", 93480, 353, 86513, 5291, 0, 2, 1676),
    @("Etiopia", 89860, 0, 43149, 45346, 0, 0, 1365),
    @("Honduras", 89381, 956, 35398, 51407, 0, 8, 2576),
    @("Bielorrusia", 88909, 619, 80503, 7469, 0, 4, 937),
    @("Venezuela", 87161, 0, 80316, 6104, 0, 0, 741),
    @("Suiza", 86167, 3008, 54600, 29422, 0, 7, 2145),
    @("China", 85704, 19, 80812, 258, 0, 0, 4634),
    @("Barein", 78224, 0, 74683, 3237, 0, 2, 304),
    @("Moldavia", 67958, 656, 49083, 17258, 0, 17, 1617),
    @("Austria", 67451, 1524, 51407, 15130, 0, 10, 914),
    @("Armenia", 66694, 1234, 48734, 16859, 0, 10, 1101),
    @("Uzbekistan", 63831, 308, 60910, 2387, 0, 3, 534),
    @("Libano", 62944, 0, 28855, 33563, 0, 0, 526),
    @("Nigeria", 61558, 0, 56697, 3736, 0, 0, 1125),
    @("Singapur", 57921, 6, 57819, 74, 0, 0, 28),
    @("Paraguay", 55452, 0, 36663, 17582, 0, 0, 1207),
    @("Argelia", 54616, 0, 38215, 14536, 0, 0, 1865),
    @("Kirguistan", 52910, 384, 46139, 5658, 0, 2, 1113),
    @("Irlanda", 52256, 1263, 23364, 27027, 0, 13, 1865),
    @("Libia", 50906, 957, 27832, 22328, 0, 14, 746),
    @("Hungria", 48757, 989, 14637, 32909, 0, 38, 1211),
    @("Estado de Palestina", 48129, 513, 41455, 6253, 0, 8, 421),
    @("Ghana", 47372, 0, 46664, 398, 0, 0, 310),
    @("Azerbaiyan", 45879, 584, 40272, 4972, 0, 5, 635),
    @("Kenia", 45647, 571, 32522, 12283, 0, 3, 842),
    @("Tunez", 42727, 0, 5032, 37008, 0, 0, 687),
    @("Jordania", 40972, 2035, 7123, 33435, 0, 34, 414),
    @("Afganistan", 40357, 70, 33790, 5068, 0, 2, 1499),
    @("Birmania", 38502, 1297, 18874, 18683, 0, 31, 945),
    @("Serbia", 36608, 326, 31536, 4292, 0, 2, 780),
    @("Dinamarca", 36373, 529, 29998, 5687, 0, 2, 688),
    @("Bosnia y Herzegovina", 35389, 728, 25560, 8812, 0, 20, 1017),
    @("El Salvador", 31975, 149, 27453, 3593, 0, 3, 929),
    @("Eslovaquia", 31400, 705, 8004, 23298, 0, 6, 98),
    @("Bulgaria", 30527, 0, 17153, 12366, 0, 0, 1008),
    @("Australia", 27429, 30, 25126, 1398, 0, 0, 905),
    @("Croacia", 26863, 890, 20962, 5519, 0, 8, 382),
    @("Grecia", 26469, 667, 9989, 15952, 0, 8, 528),
    @("Corea del Sur", 25333, 58, 23466, 1420, 0, 3, 447),
    @("Republica de Macedonia", 24196, 408, 17673, 5673, 0, 4, 850),
    @("Malasia", 22225, 862, 14351, 7681, 0, 3, 193),
    @("Camerun", 21506, 0, 20117, 965, 0, 0, 424),
    @("Costa de Marfil", 20324, 0, 20029, 174, 0, 0, 121),
    @("Georgia", 19857, 1194, 8666, 11033, 0, 15, 158),
    @("Albania", 17651, 301, 10225, 6968, 0, 4, 458),
    @("Madagascar", 16810, 0, 16215, 357, 0, 0, 238),
    @("Noruega", 16679, 76, 11863, 4538, 0, 0, 278),
    @("Zambia", 15982, 85, 15038, 598, 0, 0, 346),
    @("Montenegro", 15892, 132, 11581, 4064, 0, 7, 247),
    @("Senegal", 15459, 27, 13922, 1217, 0, 1, 320),
    @("Eslovenia", 14473, 794, 6572, 7709, 0, 2, 192),
    @("Finlandia", 13849, 294, 9100, 4398, 0, 0, 351),
    @("Sudan", 13724, 0, 6764, 6124, 0, 0, 836),
    @("Namibia", 12326, 0, 10426, 1769, 0, 0, 131),
    @("Guinea", 11538, 20, 10447, 1021, 0, 0, 70),
    @("Luxemburgo", 11241, 231, 8471, 2634, 0, 1, 136),
    @("Maldivas", 11232, 0, 10201, 994, 0, 0, 37),
    @("Mozambique", 11190, 110, 9007, 2105, 0, 3, 78),
    @("Consejo Danes para los Refugiados", 11052, 0, 10357, 392, 0, 0, 303),
    @("Uganda", 10788, 97, 7066, 3625, 0, 0, 97),
    @("Tayikistan", 10574, 41, 9617, 877, 0, 0, 80),
    @("Guayana Francesa", 10268, 0, 9995, 204, 0, 0, 69),
    @("Haiti", 8976, 0, 7303, 1442, 0, 0, 231),
    @("Gabon", 8884, 0, 8452, 378, 0, 0, 54),
    @("Jamaica", 8374, 53, 4002, 4198, 0, 1, 174),
    @("Zimbabue", 8159, 0, 7683, 244, 0, 0, 232),
    @("Lituania", 7928, 202, 3276, 4534, 0, 5, 118),
    @("Angola", 7829, 0, 3031, 4550, 0, 0, 248),
    @("Cabo Verde", 7800, 0, 6620, 1093, 0, 0, 87),
    @("Mauritania", 7621, 0, 7355, 103, 0, 0, 163),
    @("Guadalupe", 7122, 0, 2199, 4827, 0, 0, 96),
    @("Cuba", 6305, 47, 5806, 372, 0, 0, 127),
    @("Bahamas", 5923, 150, 3525, 2274, 0, 1, 124),
    @("Malaui", 5861, 1, 4757, 922, 0, 1, 182),
    @("Suazilandia", 5788, 0, 5427, 245, 0, 0, 116),
    @("Sri Lanka", 5685, 60, 3457, 2215, 0, 0, 13),
    @("Botsuana", 5609, 0, 915, 4673, 0, 0, 21),
    @("Republica de Yibuti", 5469, 0, 5379, 29, 0, 0, 61),
    @("Nicaragua", 5353, 0, 4225, 974, 0, 0, 154),
    @("Trinidad yTobago", 5298, 0, 3696, 1505, 0, 0, 97),
    @("Hong Kong", 5262, 5, 4996, 161, 0, 0, 105),
    @("Congo", 5156, 0, 3887, 1177, 0, 0, 92),
    @("Siria", 5134, 0, 1565, 3318, 0, 0, 251),
    @("Surinam", 5133, 0, 4959, 65, 0, 0, 109),
    @("Guinea Ecuatorial", 5074, 4, 4954, 37, 0, 0, 83),
    @("Ruanda", 4992, 0, 4797, 161, 0, 0, 34),
    @("Reunion", 4921, 0, 4445, 459, 0, 0, 17),
    @("Malta", 4871, 134, 3282, 1543, 0, 1, 46),
    @("Republica de Africa Central", 4856, 0, 1924, 2870, 0, 0, 62),
    @("Polinesia Francesa", 4548, 0, 3202, 1330, 0, 0, 16),
    @("Aruba", 4334, 0, 4040, 260, 0, 0, 34),
    @("Islandia", 4193, 92, 2930, 1252, 0, 0, 11),
    @("Mayotte", 4159, 0, 2964, 1152, 0, 0, 43),
    @("Estonia", 4127, 42, 3270, 786, 0, 3, 71),
    @("Somalia", 3890, 0, 3089, 702, 0, 0, 99),
    @("Guyana", 3765, 0, 2749, 905, 0, 0, 111),
    @("Tailandia", 3700, 9, 3491, 150, 0, 0, 59),
    @("Gambia", 3649, 0, 2649, 882, 0, 0, 118),
    @("Principado de Andorra", 3623, 0, 2273, 1288, 0, 0, 62),
    @("Letonia", 3609, 115, 1341, 2221, 0, 3, 47),
    @("Mali", 3407, 0, 2588, 687, 0, 0, 132),
    @("Sudan del Sur", 2847, 0, 1290, 1502, 0, 0, 55),
    @("Belice", 2833, 20, 1692, 1096, 0, 1, 45),
    @("Republica de Chipre", 2687, 0, 1444, 1218, 0, 0, 25),
    @("Uruguay", 2560, 0, 2121, 388, 0, 0, 51),
    @("Benin", 2496, 0, 2330, 125, 0, 0, 41),
    @("Guinea-Bisau", 2403, 0, 1818, 544, 0, 0, 41),
    @("Burkina Faso", 2387, 0, 1802, 520, 0, 0, 65),
    @("Sierra Leona", 2331, 0, 1760, 498, 0, 0, 73),
    @("Martinica", 2257, 0, 98, 2135, 0, 0, 24),
    @("Togo", 2071, 0, 1541, 479, 0, 0, 51),
    @("Yemen", 2056, 0, 1338, 121, 0, 0, 597),
    @("Lesoto", 1918, 85, 961, 914, 0, 1, 43),
    @("Nueva Zelanda", 1887, 1, 1829, 33, 0, 0, 25),
    @("Republica del Chad", 1390, 0, 1194, 103, 0, 0, 93),
    @("Liberia", 1381, 0, 1271, 28, 0, 0, 82),
    @("Niger", 1211, 0, 1128, 14, 0, 0, 69),
    @("Vietnam", 1141, 1, 1046, 60, 0, 0, 35),
    @("Santo Tome y Principe", 933, 0, 898, 20, 0, 0, 15),
    @("San Marino", 766, 0, 688, 36, 0, 0, 42),
    @("Curazao", 759, 8, 466, 292, 0, 0, 1),
    @("San Martin (Parte Holandesa)", 756, 0, 675, 59, 0, 0, 22),
    @("Crucero", 712, 0, 659, 40, 0, 0, 13),
    @("Islas Turcas y Caicos", 698, 0, 689, 3, 0, 0, 6),
    @("Gibraltar", 608, 31, 471, 137, 0, 0, 0),
    @("Papua Nueva Guinea", 581, 0, 541, 33, 0, 0, 7),
    @("Burundi", 550, 1, 497, 52, 0, 0, 1),
    @("Taiwan", 543, 3, 495, 41, 0, 0, 7),
    @("San Martin (Parte Francesa)", 531, 0, 380, 143, 0, 0, 8),
    @("Tanzania", 509, 0, 183, 305, 0, 0, 21),
    @("Comoras", 502, 0, 485, 10, 0, 0, 7),
    @("Islas Feroe", 488, 3, 472, 16, 0, 0, 0),
    @("Eritrea", 452, 0, 388, 64, 0, 0, 0),
    @("Mauricio", 419, 0, 379, 30, 0, 0, 10),
    @("Isla de Man", 348, 0, 321, 3, 0, 0, 24),
    @("Butan", 330, 3, 301, 29, 0, 0, 0),
    @("Mongolia", 326, 2, 312, 14, 0, 0, 0),
    @("Camboya", 285, 2, 280, 5, 0, 0, 0),
    @("Monaco", 268, 0, 222, 44, 0, 0, 2),
    @("Liechtenstein", 235, 11, 142, 92, 0, 0, 1),
    @("Islas Caimanes", 235, 0, 212, 22, 0, 0, 1),
    @("Barbados", 222, 0, 203, 12, 0, 0, 7),
    @("Bermudas", 185, 0, 172, 4, 0, 0, 9),
    @("Bonaire, San Eustaquio y Saba", 150, 0, 121, 26, 0, 0, 3),
    @("Seychelles", 149, 0, 148, 1, 0, 0, 0),
    @("Brunei", 147, 0, 143, 1, 0, 0, 3),
    @("Antigua y Barbuda", 119, 0, 101, 15, 0, 0, 3),
    @("San Bartolome", 72, 0, 55, 17, 0, 0, 0),
    @("Islas Virgenes Britanicas", 71, 0, 70, 0, 0, 0, 1),
    @("San Vicente y las Granadinas", 67, 0, 64, 3, 0, 0, 0),
    @("Macao", 46, 0, 46, 0, 0, 0, 0),
    @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
    @("Santa Lucia", 36, 0, 27, 9, 0, 0, 0),
    @("Dominica", 33, 0, 29, 4, 0, 0, 0),
    @("Guam", 32, 0, 0, 31, 0, 0, 1),
    @("Fiyi", 32, 0, 30, 0, 0, 0, 2),
    @("Timor Oriental", 29, 0, 28, 1, 0, 0, 0),
    @("Santa Sede", 27, 0, 15, 12, 0, 0, 0),
    @("Granada", 27, 0, 24, 3, 0, 0, 0),
    @("Nueva Caledonia", 27, 0, 27, 0, 0, 0, 0),
    @("Laos", 23, 0, 22, 1, 0, 0, 0),
    @("San Cristobal y Nieves", 19, 0, 19, 0, 0, 0, 0),
    @("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
    @("San Pedro y Miquelon", 16, 0, 12, 4, 0, 0, 0),
    @("Groenlandia", 16, 0, 14, 2, 0, 0, 0),
    @("Montserrat", 13, 0, 12, 0, 0, 0, 1),
    @("Islas Malvinas", 13, 0, 13, 0, 0, 0, 0),
    @("Sahara Occidental", 10, 0, 8, 1, 0, 0, 1),
    @("Islas Salomon", 3, 0, 0, 3, 0, 0, 0),
    @("Anguila", 3, 0, 3, 0, 0, 0, 0),
    @("Wallis y Futuna", 1, 0, 0, 1, 0, 0, 0)
)

$startRow = 4
for ($i = 0; $i -lt $countryData.Count; $i++) {
    $row = $startRow + $i
    $rec = $countryData[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = [double]$rec[1]
    $ws.Cells.Item($row, 3).Value = [double]$rec[2]
    $ws.Cells.Item($row, 4).Value = [double]$rec[3]
    $ws.Cells.Item($row, 5).Value = [double]$rec[4]
    $ws.Cells.Item($row, 6).Value = [double]$rec[5]
    $ws.Cells.Item($row, 7).Value = [double]$rec[6]
    $ws.Cells.Item($row, 8).Value = [double]$rec[7]
}